$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.458.61"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.838.05"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.32"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5328"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +2.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3011"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -6.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06865"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.43"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -7.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.870.09"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7339"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -5.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07442"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -4.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.21"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.967"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.91"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007906"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.498.51"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.584"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.972"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.246"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -2.31%  "
$ws.Range("B24").Value = "LidoDAOToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.219"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.52"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.683"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.91"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "110.30"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.250"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08778"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.021"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04785"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.922"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7292"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.127"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.098"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.288"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01709"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -4.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.4701"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9044"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.64"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -3.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.868"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.348"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -3.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.026"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4079"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -3.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1231"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.81"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05800"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8922"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.10"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.23%  "
